$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 (shifts existing rows 10-34 down to 11-35)
$ws.Rows("10").Insert()

# Populate the new row 10 as a duplicate of the "Sword World RPG World Guide" entry
# (row 9) but pointing at the alternate cover image.
$ws.Range("A10").Value = 1993
$ws.Range("B10").Value = "ソード・ワールドRPGワールドガイド"
$ws.Range("C10").Value = "Sword World RPG World Guide"
$ws.Range("D10").Value = "Fujimi Shobo"
$ws.Range("E10").Value = "sword_world_world_guide_alt.jpg"

# Match the author's final selection position
$ws.Range("E11").Select()
